$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (shifts old J,K -> K,L)
$ws.Range("J1:J2").EntireColumn.Insert()

# Populate the new column J with the tenant_id comment (row1) and value template (row2)
$ws.Range("J1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'
$ws.Range("J2").Value = '<%=model.tenant_id_lbl%>'
